# Updated with WorkFlow Execution
# The "StoreResponseVariables" column (L) for the PetGet test case (row 3)
# gets an additional stored-variable mapping for the category name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").Value = "petId=id;petName=name;category_name=category.name"
